$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2-11 hold a "Förändrad" date that was updated from
# 2023-10-13 (45212) to 2023-10-22 (45221) for all current data rows.
$newDate = Get-Date -Year 2023 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
